$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Valor Depois" values for rows 7 and 8 (spreadsheet rows)
$ws.Range("E7").Value = 6
$ws.Range("E8").Value = 23

# Update the active selection to E8 (matches the recorded cursor position after edit)
$ws.Range("E8").Select()
